# Update on 2019-03-06, 支出生活费400+1200
# Target worksheet is "第二学年" (second academic year), the 2nd sheet in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Row 27: new expense entry #24 - 生活费 (living expenses) 400, dated 2019-02-20
$ws.Cells.Item(27, 2).Value = 24
$ws.Cells.Item(27, 3).Value = "支出"
$ws.Cells.Item(27, 4).Value = 400
$ws.Cells.Item(27, 5).Value = 43516
$ws.Cells.Item(27, 6).Value = "生活费"
$ws.Cells.Item(27, 7).Value = "生活费(2019-02-20 到 2019-02-28)"

# Row 28: new expense entry #25 - 生活费 (living expenses) 1200, dated 2019-03-03
$ws.Cells.Item(28, 2).Value = 25
$ws.Cells.Item(28, 3).Value = "支出"
$ws.Cells.Item(28, 4).Value = 1200
$ws.Cells.Item(28, 5).Value = 43527
$ws.Cells.Item(28, 6).Value = "生活费"
$ws.Cells.Item(28, 7).Value = "生活费(2019-03-01 到 2019-03-31)"

# Recalculate the workbook so SUMIFS totals (I4/J4/K4, J10, etc.) reflect the new rows
$excel.Calculate()

# Move the active selection to G29, matching the author's cursor position after the edit
$ws.Range("G29").Select()
